$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 23
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12)
$ws.Range("C2:C23").Value = 45181
